# Insert a new weekly price-observation row for Albahaca ("Terminal La
# Palmera de La Serena") just above the current row 55, pushing all the
# rows below it (old rows 55-110) down by one. The used range grows from
# A1:R110 to A1:R111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 55..110 down to 56..111 by inserting a fresh row at 55.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new observation. The
# non-numeric / non-date columns mirror the row that used to sit at 55
# (same market, region, product, quality, prices bracket, unit and
# origin) - only the date (D) and volume (J) are new values.
$ws.Cells.Item(55, 1).Value = 8
$ws.Cells.Item(55, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(55, 3).Value = "Coquimbo"
$ws.Cells.Item(55, 4).Value = 44790
$ws.Cells.Item(55, 5).Value = 4
$ws.Cells.Item(55, 6).Value = 100112052
$ws.Cells.Item(55, 7).Value = "Albahaca"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 1600
$ws.Cells.Item(55, 11).Value = 3300
$ws.Cells.Item(55, 12).Value = 3500
$ws.Cells.Item(55, 13).Value = 3400
$ws.Cells.Item(55, 14).Value = "`$/paquete"
$ws.Cells.Item(55, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value = 3400
$ws.Cells.Item(55, 17).Value = 1
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat
